$wb = $excel.ActiveWorkbook

$wsCheckReport = $wb.Worksheets.Item("checkReport")
$wsViewMore = $wb.Worksheets.Item("viewMore")

$newDate = "2023-07-07"

# Row 2 - new "Leather00-13" item duplicated into the report.
$wsCheckReport.Cells.Item(2, 5).Value = $newDate
$wsCheckReport.Cells.Item(2, 1).Value = "Leather00-13"
$wsCheckReport.Cells.Item(2, 2).Value = "LT-013"
$wsCheckReport.Cells.Item(2, 3).Value = "Mulchand & Sons"

# Rows 3-7 - automation items shifted from 65-69/64 up to 70-74, plus a new "Vendor" value.
$itemNames   = @("(Automation) Item 70", "(Automation) Item 71", "(Automation) Item 72", "(Automation) Item 73", "(Automation) Item 74")
$productCodes = @("AT-IT-70", "AT-IT-71", "AT-IT-72", "AT-IT-73", "AT-IT-74")
$vendors     = @("AKM Leather", "Mulchand & Sons", "Tamang Suppliers", "Silver Zipper Supplier", "Vendor")

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 3
    $wsCheckReport.Cells.Item($row, 1).Value = $itemNames[$i]
    $wsCheckReport.Cells.Item($row, 2).Value = $productCodes[$i]
    $wsCheckReport.Cells.Item($row, 3).Value = $vendors[$i]
    $wsCheckReport.Cells.Item($row, 5).Value = $newDate
}

# Duplicate the same data (rows 2-7) into the "viewMore" sheet.
for ($row = 2; $row -le 7; $row++) {
    $wsViewMore.Cells.Item($row, 1).Value = $wsCheckReport.Cells.Item($row, 1).Value2
    $wsViewMore.Cells.Item($row, 2).Value = $wsCheckReport.Cells.Item($row, 2).Value2
    $wsViewMore.Cells.Item($row, 3).Value = $wsCheckReport.Cells.Item($row, 3).Value2
    $wsViewMore.Cells.Item($row, 5).Value = $wsCheckReport.Cells.Item($row, 5).Value2
}

# Update the selection on both sheets to C7, and make "checkReport" the active tab.
$wsViewMore.Range("C7").Select()
$wsCheckReport.Activate()
$wsCheckReport.Range("C7").Select()
